$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") should carry the same style as
# the existing header cells (bold font, border, centered/top alignment).
# Copy H1 -> I1 and H1 -> J1 first to pick up that formatting, then set
# the actual header text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells I2 and J2 are plain numeric values (no special style).
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
